# CronogramaFinal.xlsx - layout adjustment after meeting with Ricardo.
# - Recolor the Gantt categories with a new green palette.
# - Insert a new "Resultados" entry (2021-11-01 .. 2021-11-30).
# - Move selection to I24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row just before the old row 23 ("Revisao" block) so the
#        new "Resultados" data point lands at row 23, pushing the rest down.
$ws.Rows(23).Insert()

$ws.Range("A23").Value2 = "Resultados"
$ws.Range("B23").Value2 = "Resultados"
$ws.Range("C23").Value2 = 44501
$ws.Range("D23").Value2 = 44530

# --- 2. New color palette, keyed by the task label (column A/B).
#        Every row with a given label gets the same fill/text color pair.
$colorMap = @{
    "Análise "        = @("#4B6043", "white")
    "Coleta "         = @("#DDEAD1", "black")
    "Interpretação "  = @("#658354", "white")
    "Resultados"      = @("#C7DDB5", "black")
    "Revisão"         = @("#95BB72", "white")
    "Sessão Pública " = @("#B3CF99", "black")
}

for ($r = 2; $r -le 30; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    $colors = $colorMap[$label]
    if ($colors) {
        $ws.Cells.Item($r, 5).Value2 = $colors[0]
        $ws.Cells.Item($r, 6).Value2 = $colors[1]
    }
}

# --- 3. Selection as left by the editing session.
$null = $ws.Range("I24").Select()
